# Updated symbol list on Tue Jan  3 16:54:37 UTC 2023 with GitHub Actions
# Refresh the scraped coin price/volume data (and the coin order for rows
# 10-15) to match the latest run of the scraper.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}


Set-TextValue "D2" '245.15'
Set-TextValue "E2" '-0.68%'
Set-TextValue "D3" '28.86'
Set-TextValue "E3" '-2.09%'
Set-TextValue "D4" '5.286'
Set-TextValue "E4" '2.25%'
Set-TextValue "D5" '0.05714'
Set-TextValue "E5" '0.19%'
Set-TextValue "E6" '0.33%'
Set-TextValue "D7" '3.183'
Set-TextValue "E7" '4.00%'
Set-TextValue "D8" '0.8533'
Set-TextValue "E8" '-0.61%'
Set-TextValue "D9" '0.8582'
Set-TextValue "E9" '-1.56%'
Set-TextValue "B10" 'One'
Set-TextValue "C10" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue "D10" '0.01007'
Set-TextValue "E10" '1,575.11%'
Set-TextValue "B11" 'WazirX'
Set-TextValue "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D11" '0.1371'
Set-TextValue "E11" '0.33%'
Set-TextValue "B12" 'MandalaExchangeToken'
Set-TextValue "C12" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D12" '0.07039'
Set-TextValue "E12" '-0.47%'
Set-TextValue "B13" 'BitrueCoin'
Set-TextValue "C13" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D13" '0.03159'
Set-TextValue "E13" '9.96%'
Set-TextValue "B14" 'BitMartToken'
Set-TextValue "C14" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D14" '0.09302'
Set-TextValue "E14" '-0.84%'
Set-TextValue "B15" 'BitForexToken'
Set-TextValue "C15" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D15" '0.001532'
Set-TextValue "E15" '1.44%'
Set-TextValue "D16" '0.006049'
Set-TextValue "E16" '-0.28%'
Set-TextValue "D17" '3.492'
Set-TextValue "E17" '0.22%'
Set-TextValue "D18" '2.174'
Set-TextValue "E18" '-4.66%'
Set-TextValue "D19" '0.3164'
Set-TextValue "E19" '-0.25%'
Set-TextValue "D20" '0.03345'
Set-TextValue "E20" '1.23%'
Set-TextValue "E21" '-1.78%'
Set-TextValue "D22" '3.502'
Set-TextValue "E22" '0.91%'
Set-TextValue "D23" '0.04089'
Set-TextValue "E23" '-2.08%'
Set-TextValue "D24" '0.1380'
Set-TextValue "E24" '0.02%'
Set-TextValue "D25" '0.001225'
Set-TextValue "E25" '0.50%'
Set-TextValue "D26" '0.004143'
Set-TextValue "E26" '-17.73%'
Set-TextValue "E27" '-0.77%'
Set-TextValue "D28" '0.0001450'
Set-TextValue "E28" '-25.22%'
Set-TextValue "D40" '0.03764'
Set-TextValue "E40" '0.26%'
Set-TextValue "D41" '0.1064'
Set-TextValue "E41" '-0.70%'
Set-TextValue "D42" '0.003703'
Set-TextValue "E42" '-36.14%'
Set-TextValue "D43" '0.002450'
Set-TextValue "E43" '16.67%'
Set-TextValue "D44" '0.009348'
Set-TextValue "E44" '-8.47%'
Set-TextValue "D45" '0.00005311'
Set-TextValue "E45" '2.73%'
Set-TextValue "E46" '0.07%'
Set-TextValue "D47" '0.07504'
Set-TextValue "E47" '7.21%'
Set-TextValue "D48" '0.002442'
Set-TextValue "E48" '-4.92%'
Set-TextValue "D49" '0.00002101'
Set-TextValue "E49" '0.07%'
Set-TextValue "D50" '0.0002001'
Set-TextValue "E50" '0.07%'
